$d = $word.ActiveDocument

# Merge the multiple runs that make up the dedication line into a single
# run by replacing the full text with itself. Word's Find/Replace collapses
# the matched range into one run (using the formatting of the first run)
# and drops any bookmarks that fall inside the replaced range.
$d.Content.Find.Execute(
    "To my family: Johanna, Kaitlyn, Chris, and Amy",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "To my family: Johanna, Kaitlyn, Chris, and Amy",
    2
)
